$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.394.43"
$ws.Range("E2").Value = "  +3.28%  "
$ws.Range("D3").Value = "1.870.46"
$ws.Range("E3").Value = "  +1.54%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'339.46"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.4704"
$ws.Range("E7").Value = "  +1.86%  "
$ws.Range("D8").Value = "'0.3958"
$ws.Range("E8").Value = "  +3.31%  "
$ws.Range("D9").Value = "'47.47"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").Value = "'0.08016"
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("D12").Value = "'21.84"
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").Value = "1.871.37"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").Value = "'7.239"
$ws.Range("E15").Value = "  +2.98%  "
$ws.Range("D16").Value = "'91.28"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").Value = "'0.00001039"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").Value = "'0.06619"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "'17.54"
$ws.Range("E20").Value = "  +3.47%  "
$ws.Range("E21").Value = "  -0.33%  "
$ws.Range("D22").Value = "28.398.37"
$ws.Range("E22").Value = "  +3.28%  "
$ws.Range("E23").Value = "  +2.02%  "
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").Value = "2.081.24"
$ws.Range("E26").Value = "  +1.35%  "
$ws.Range("D27").Value = "'160.29"
$ws.Range("E27").Value = "  +2.02%  "
$ws.Range("D28").Value = "'19.74"
$ws.Range("E28").Value = "  +1.84%  "
$ws.Range("D29").Value = "'2.131"
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "'5.516"
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("D31").Value = "'119.99"
$ws.Range("E31").Value = "  +0.90%  "
$ws.Range("D32").Value = "'0.9682"
$ws.Range("E32").Value = "  +1.17%  "
$ws.Range("D33").Value = "'0.09481"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "'3.573"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "'1.378"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("E36").Value = "  +2.03%  "
$ws.Range("D37").Value = "'0.06080"
$ws.Range("E37").Value = "  +2.33%  "
$ws.Range("D38").Value = "'0.02246"
$ws.Range("D39").Value = "'8.375"
$ws.Range("E39").Value = "  +3.37%  "
$ws.Range("E40").Value = "  +2.41%  "
$ws.Range("D41").Value = "'0.5944"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("E43").Value = "  +1.45%  "
$ws.Range("E44").Value = "  +2.25%  "
$ws.Range("D45").Value = "'1.293"
$ws.Range("E45").Value = "  +3.61%  "
$ws.Range("D46").Value = "'0.5575"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("E47").Value = "  +1.70%  "
$ws.Range("D48").Value = "'1.953"
$ws.Range("E48").Value = "  +4.44%  "
$ws.Range("D49").Value = "'0.06854"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").Value = "'2.048"
$ws.Range("E50").Value = "  +14.67%  "
$ws.Range("D51").Value = "'111.33"
$ws.Range("E51").Value = "  +0.97%  "
